$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-06-21 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-22 Thursday", 2) | Out-Null

# Update each math-problem cell in the table (row-major order, 5 columns)
$t = $d.Tables.Item(1)

$r = $t.Cell(1,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "1+89="
$r = $t.Cell(1,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "57-2="
$r = $t.Cell(1,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "12+29="
$r = $t.Cell(1,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "27+18="
$r = $t.Cell(1,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "66-63="
$r = $t.Cell(2,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "23+54="
$r = $t.Cell(2,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "19+17="
$r = $t.Cell(2,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "54+4="
$r = $t.Cell(2,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "6+11="
$r = $t.Cell(2,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "55-41="
$r = $t.Cell(3,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "44+5="
$r = $t.Cell(3,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "76-27="
$r = $t.Cell(3,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "78-14="
$r = $t.Cell(3,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "69-54="
$r = $t.Cell(3,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "41+52="
$r = $t.Cell(4,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "28-1="
$r = $t.Cell(4,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "92-88="
$r = $t.Cell(4,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "34-8="
$r = $t.Cell(4,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "88-9="
$r = $t.Cell(4,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "49-26="
$r = $t.Cell(5,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "51+12="
$r = $t.Cell(5,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "88-84="
$r = $t.Cell(5,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "94-82="
$r = $t.Cell(5,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "37+45="
$r = $t.Cell(5,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "57-14="
$r = $t.Cell(6,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "24+5="
$r = $t.Cell(6,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "24-9="
$r = $t.Cell(6,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "92-60="
$r = $t.Cell(6,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "36+28="
$r = $t.Cell(6,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "51+24="
$r = $t.Cell(7,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "41-16="
$r = $t.Cell(7,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "82-80="
$r = $t.Cell(7,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "93-45="
$r = $t.Cell(7,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "11+11="
$r = $t.Cell(7,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "64+18="
$r = $t.Cell(8,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "21-12="
$r = $t.Cell(8,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "36-31="
$r = $t.Cell(8,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "48-43="
$r = $t.Cell(8,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "16+77="
$r = $t.Cell(8,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "18+53="
$r = $t.Cell(9,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "4+30="
$r = $t.Cell(9,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "87-70="
$r = $t.Cell(9,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "98-18="
$r = $t.Cell(9,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "71-55="
$r = $t.Cell(9,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "59+12="
$r = $t.Cell(10,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "98-66="
$r = $t.Cell(10,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "4+81="
$r = $t.Cell(10,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "96-41="
$r = $t.Cell(10,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "10+46="
$r = $t.Cell(10,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "73-47="
$r = $t.Cell(11,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "14+51="
$r = $t.Cell(11,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "67+17="
$r = $t.Cell(11,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "68-19="
$r = $t.Cell(11,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "53-43="
$r = $t.Cell(11,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "58-50="
$r = $t.Cell(12,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "60-53="
$r = $t.Cell(12,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "44-20="
$r = $t.Cell(12,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "94+4="
$r = $t.Cell(12,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "22+10="
$r = $t.Cell(12,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "95-48="
$r = $t.Cell(13,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "46-31="
$r = $t.Cell(13,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "3+88="
$r = $t.Cell(13,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "21+29="
$r = $t.Cell(13,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "49-44="
$r = $t.Cell(13,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "61+5="
$r = $t.Cell(14,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "93-72="
$r = $t.Cell(14,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "52-18="
$r = $t.Cell(14,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "99-96="
$r = $t.Cell(14,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "20+72="
$r = $t.Cell(14,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "53-29="
$r = $t.Cell(15,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "53-18="
$r = $t.Cell(15,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "66+6="
$r = $t.Cell(15,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "30+39="
$r = $t.Cell(15,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "6+84="
$r = $t.Cell(15,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "41-32="
$r = $t.Cell(16,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "75-53="
$r = $t.Cell(16,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "7+63="
$r = $t.Cell(16,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "46+2="
$r = $t.Cell(16,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "49+20="
$r = $t.Cell(16,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "6+48="
$r = $t.Cell(17,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "87-55="
$r = $t.Cell(17,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "42+12="
$r = $t.Cell(17,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "43-5="
$r = $t.Cell(17,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "72+17="
$r = $t.Cell(17,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "62-53="
$r = $t.Cell(18,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "20+60="
$r = $t.Cell(18,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "76+20="
$r = $t.Cell(18,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "30+49="
$r = $t.Cell(18,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "97-1="
$r = $t.Cell(18,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "6+49="
$r = $t.Cell(19,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "96-60="
$r = $t.Cell(19,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "47+17="
$r = $t.Cell(19,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "50+44="
$r = $t.Cell(19,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "68+25="
$r = $t.Cell(19,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "49-44="
$r = $t.Cell(20,1).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "8+85="
$r = $t.Cell(20,2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "74-24="
$r = $t.Cell(20,3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "90-59="
$r = $t.Cell(20,4).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "72-70="
$r = $t.Cell(20,5).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "98-79="

Write-Host "Replacements complete"
